# actualizacion libretas bancarias modal v11
# Adds two new movement rows (30 and 31) to the bank-statement sheet,
# copying the formatting of the last existing row (29) and filling in
# the new values, then updates the sheet view to reflect where the
# user ended up (scrolled down, with the newly added cell selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 : 02-07-2020 / PRUEBA TEST 2 -------------------------------
$ws.Rows.Item(29).Copy()
$ws.Rows.Item(30).Insert()
$ws.Range("A30").Value = 44014
$ws.Range("B30").Value = "PRUEBA TEST 2"
$ws.Range("C30").Value = 41438426
$ws.Range("D30").Value = 100
$ws.Range("E30").Value = "TJA"

# --- Row 31 : 04-07-2020 / PRUEBA TEST 4 -------------------------------
$ws.Rows.Item(30).Copy()
$ws.Rows.Item(31).Insert()
$ws.Range("A31").Value = 44016
$ws.Range("B31").Value = "PRUEBA TEST 4"
$ws.Range("C31").Value = 41438426
$ws.Range("D31").Value = 100
$ws.Range("E31").Value = "TJA"

$excel.Application.CutCopyMode = $false

# Leave the view scrolled to the new rows with the last entry selected,
# matching where the user was working in the workbook.
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 17
